$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.203.76"
$ws.Range("E2").Value = "  +0.80%  "
$ws.Range("D3").Value = "1.803.91"
$ws.Range("E3").Value = "  +3.52%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "336.33"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "0.4667"
$ws.Range("E7").Value = "  +24.01%  "
$ws.Range("D8").Value = "0.3713"
$ws.Range("E8").Value = "  +10.93%  "
$ws.Range("D9").Value = "45.46"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.07676"
$ws.Range("E10").Value = "  +6.54%  "
$ws.Range("D11").Value = "1.156"
$ws.Range("E11").Value = "  +3.99%  "
$ws.Range("D12").Value = "22.65"
$ws.Range("E12").Value = "  +1.18%  "
$ws.Range("D13").Value = "1.002"
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "6.379"
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("E15").Value = "  +4.23%  "
$ws.Range("D16").Value = "1.795.05"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("E17").Value = "  +3.94%  "
$ws.Range("D18").Value = "0.06729"
$ws.Range("E18").Value = "  +2.41%  "
$ws.Range("D19").Value = "82.77"
$ws.Range("E19").Value = "  +3.86%  "
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("D21").Value = "17.47"
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("D22").Value = "6.433"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").Value = "28.138.84"
$ws.Range("E23").Value = "  +0.54%  "
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").Value = "2.418"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").Value = "20.83"
$ws.Range("E26").Value = "  +5.31%  "
$ws.Range("D27").NumberFormat = "@"  # preserve trailing zero(s) as text
$ws.Range("D27").Value = "2.400"
$ws.Range("E27").Value = "  +3.44%  "
$ws.Range("D28").Value = "152.38"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "2.006.75"
$ws.Range("E29").Value = "  +2.94%  "
$ws.Range("D30").NumberFormat = "@"  # preserve trailing zero(s) as text
$ws.Range("D30").Value = "134.70"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").NumberFormat = "@"  # preserve trailing zero(s) as text
$ws.Range("D31").Value = "1.270"
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"  # preserve trailing zero(s) as text
$ws.Range("D33").Value = "0.09660"
$ws.Range("E33").Value = "  +10.73%  "
$ws.Range("D34").Value = "5.922"
$ws.Range("E34").Value = "  +2.54%  "
$ws.Range("D35").Value = "0.2248"
$ws.Range("E35").Value = "  +6.73%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").Value = "0.02377"
$ws.Range("E37").Value = "  +2.76%  "
$ws.Range("D38").Value = "0.06405"
$ws.Range("E38").Value = "  +3.38%  "
$ws.Range("D39").Value = "0.6727"
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("D40").Value = "5.273"
$ws.Range("E40").Value = "  +2.83%  "
$ws.Range("D41").Value = "1.528"
$ws.Range("E41").Value = "  +6.22%  "
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("D43").Value = "8.151"
$ws.Range("D44").Value = "14.21"
$ws.Range("E44").Value = "  +3.79%  "
$ws.Range("D45").Value = "0.9994"
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("D46").Value = "0.6185"
$ws.Range("E46").Value = "  +2.47%  "
$ws.Range("D47").Value = "3.841"
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "130.21"
$ws.Range("E48").Value = "  +1.90%  "
$ws.Range("D49").Value = "2.068"
$ws.Range("E49").Value = "  +2.64%  "
$ws.Range("D50").Value = "1.188"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").NumberFormat = "@"  # preserve trailing zero(s) as text
$ws.Range("D51").Value = "0.07150"
$ws.Range("E51").Value = "  +0.37%  "
